# Insert a new data row at row 177 (pushing the existing rows 177-261 down
# to 178-262) and populate it with a new "Alcachofa" price record, matching
# the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(177).Insert()

$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").Value = 44839
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = 100112013
$ws.Range("G177").Value = "Alcachofa"
$ws.Range("H177").Value = "Madrigal"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 120
$ws.Range("K177").Value = 11000
$ws.Range("L177").Value = 11000
$ws.Range("M177").Value = 11000
$ws.Range("N177").Value = "`$/caja 40 unidades"
$ws.Range("O177").Value = "Región Metropolitana"
$ws.Range("P177").Value = 275
$ws.Range("Q177").Value = 40
$ws.Range("R177").Value = "Hortaliza"
